# Add a "metadata" tab after the existing "data" sheet, matching the
# panelapp export's metadata schema, and mirror the header/index-cell
# formatting used on the "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Copy the bold/bordered/centered header style from the "data" sheet
$data.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row
$ws.Range("A2").Value = 0
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Catecholaminergic Polymorphic Ventricular Tachycardia"
$ws.Range("C2").Value = 92

# Force D2 to be stored as text "0.26" rather than the number 0.26
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.26"

$ws.Range("E2").Value = "2020-06-03T21:09:12.508772Z"
$ws.Range("F2").Value = "2021-10-05 14:33:23.325730"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/92/?format=json"

# Refresh the "time_taken" values on the "data" sheet to reflect the
# re-run query timestamps.
$data.Range("F2").Value = "2021-10-05 14:33:23.328267"
$data.Range("F3").Value = "2021-10-05 14:33:23.328272"
$data.Range("F4").Value = "2021-10-05 14:33:23.328275"
$data.Range("F5").Value = "2021-10-05 14:33:23.328277"
$data.Range("F6").Value = "2021-10-05 14:33:23.328279"
$data.Range("F7").Value = "2021-10-05 14:33:23.328282"
$data.Range("F8").Value = "2021-10-05 14:33:23.328284"
$data.Range("F9").Value = "2021-10-05 14:33:23.328285"

Write-Output "metadata sheet added"
